$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the participant row for "remas ali almadani_20251202_134128"
# (row 22), which was a duplicate submission removed from the admin panel.
# All subsequent rows shift up by one.
$ws.Rows.Item(22).Delete()
